$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028528782021803
$ws.Range("D2").Value = 1.037627011731702
$ws.Range("E2").Value = 1.037150153806133
$ws.Range("F2").Value = 1.044963062975236
$ws.Range("I2").Value = 1.036838182609293
$ws.Range("J2").Value = 1.033680597415674
$ws.Range("K2").Value = 1.040417357087322
$ws.Range("L2").Value = 1.03994186036851
$ws.Range("M2").Value = 1.047732637222567
$ws.Range("N2").Value = 1.015244064075143

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029382535204035
$ws.Range("D3").Value = 1.038149895012419
$ws.Range("E3").Value = 1.037942831302007
$ws.Range("F3").Value = 1.045917564499063
$ws.Range("I3").Value = 1.037010305478713
$ws.Range("J3").Value = 1.034175488631602
$ws.Range("K3").Value = 1.040750846638599
$ws.Range("L3").Value = 1.04054433141454
$ws.Range("M3").Value = 1.048498109937444
$ws.Range("N3").Value = 1.015408258463826

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029935329201887
$ws.Range("D4").Value = 1.038487958456991
$ws.Range("E4").Value = 1.038456524174829
$ws.Range("F4").Value = 1.046536309035299
$ws.Range("I4").Value = 1.037120190943779
$ws.Range("J4").Value = 1.034495431615556
$ws.Range("K4").Value = 1.040965647602495
$ws.Range("L4").Value = 1.040934292407169
$ws.Range("M4").Value = 1.048993917284953
$ws.Range("N4").Value = 1.015514386075599

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030167808241237
$ws.Range("D5").Value = 1.038630012042195
$ws.Range("E5").Value = 1.038672665181523
$ws.Range("F5").Value = 1.046796695307102
$ws.Range("I5").Value = 1.037166029279617
$ws.Range("J5").Value = 1.0346298663287
$ws.Range("K5").Value = 1.041055711705317
$ws.Range("L5").Value = 1.041098259914119
$ws.Range("M5").Value = 1.049202471820275
$ws.Range("N5").Value = 1.01555897366583

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030206847405779
$ws.Range("D6").Value = 1.038653859402328
$ws.Range("E6").Value = 1.038708966981029
$ws.Range("F6").Value = 1.046840430862726
$ws.Range("I6").Value = 1.037173704746484
$ws.Range("J6").Value = 1.034652434410209
$ws.Range("K6").Value = 1.041070819853146
$ws.Range("L6").Value = 1.041125792378987
$ws.Range("M6").Value = 1.04923749587054
$ws.Range("N6").Value = 1.01556645843757

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029938435268148
$ws.Range("D7").Value = 1.038489856854767
$ws.Range("E7").Value = 1.038459411538526
$ws.Range("F7").Value = 1.046539787285991
$ws.Range("I7").Value = 1.037120804843802
$ws.Range("J7").Value = 1.034497228212593
$ws.Range("K7").Value = 1.040966851981311
$ws.Range("L7").Value = 1.04093648324099
$ws.Range("M7").Value = 1.048996703542458
$ws.Range("N7").Value = 1.015514981969503

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028817236810443
$ws.Range("D8").Value = 1.03780377875667
$ws.Range("E8").Value = 1.037417880912095
$ws.Range("F8").Value = 1.045285408714451
$ws.Range("I8").Value = 1.03689666020922
$ws.Range("J8").Value = 1.033847906480433
$ws.Range("K8").Value = 1.0405302653758
$ws.Range("L8").Value = 1.040145442257511
$ws.Range("M8").Value = 1.047991228878748
$ws.Range("N8").Value = 1.015299578342354

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026844350421973
$ws.Range("D9").Value = 1.036592786822446
$ws.Range("E9").Value = 1.035588597258563
$ws.Range("F9").Value = 1.043083669944983
$ws.Range("I9").Value = 1.036490320884332
$ws.Range("J9").Value = 1.03270158882821
$ws.Range("K9").Value = 1.039753432065143
$ws.Range("L9").Value = 1.0387525176692
$ws.Range("M9").Value = 1.046223311723763
$ws.Range("N9").Value = 1.014919130789967

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025531064044551
$ws.Range("D10").Value = 1.035784207840833
$ws.Range("E10").Value = 1.034373222544149
$ws.Range("F10").Value = 1.041621747678446
$ws.Range("I10").Value = 1.036211838901601
$ws.Range("J10").Value = 1.031936013225089
$ws.Range("K10").Value = 1.039230577592385
$ws.Range("L10").Value = 1.037824643489367
$ws.Range("M10").Value = 1.04504737992812
$ws.Range("N10").Value = 1.0146649327292

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024962881313058
$ws.Range("D11").Value = 1.035433812140337
$ws.Range("E11").Value = 1.033847955770299
$ws.Range("F11").Value = 1.040990139694846
$ws.Range("I11").Value = 1.03608946437876
$ws.Range("J11").Value = 1.031604201276876
$ws.Range("K11").Value = 1.03900301591817
$ws.Range("L11").Value = 1.037423055689323
$ws.Range("M11").Value = 1.044538842020059
$ws.Range("N11").Value = 1.014554732957943

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024751906258967
$ws.Range("D12").Value = 1.035303620147436
$ws.Range("E12").Value = 1.033652999929034
$ws.Range("F12").Value = 1.040755746423384
$ws.Range("I12").Value = 1.03604374098915
$ws.Range("J12").Value = 1.031480905697039
$ws.Range("K12").Value = 1.03891831619134
$ws.Range("L12").Value = 1.037273917587655
$ws.Range("M12").Value = 1.044350047231685
$ws.Range("N12").Value = 1.014513780681857

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.024797157784993
$ws.Range("D13").Value = 1.035331548524266
$ws.Range("E13").Value = 1.033694811731502
$ws.Range("F13").Value = 1.040806014863318
$ws.Range("I13").Value = 1.036053560926812
$ws.Range("J13").Value = 1.031507355074133
$ws.Range("K13").Value = 1.038936492384622
$ws.Range("L13").Value = 1.037305906862006
$ws.Range("M13").Value = 1.044390539859287
$ws.Range("N13").Value = 1.014522565945252

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024945440545908
$ws.Range("D14").Value = 1.035423051222457
$ws.Range("E14").Value = 1.033831837557195
$ws.Range("F14").Value = 1.040970760297872
$ws.Range("I14").Value = 1.036085690332017
$ws.Range("J14").Value = 1.031594010554573
$ws.Range("K14").Value = 1.038996018137113
$ws.Range("L14").Value = 1.037410727278448
$ws.Range("M14").Value = 1.044523234159513
$ws.Range("N14").Value = 1.014551348218829

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025036812203824
$ws.Range("D15").Value = 1.035479423889424
$ws.Range("E15").Value = 1.033916283830967
$ws.Range("F15").Value = 1.04107229381028
$ws.Range("I15").Value = 1.036105450827
$ws.Range("J15").Value = 1.031647395814349
$ws.Range("K15").Value = 1.039032671011114
$ws.Range("L15").Value = 1.037475314541839
$ws.Range("M15").Value = 1.044605004636514
$ws.Range("N15").Value = 1.014569079398874

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025568782667526
$ws.Range("D16").Value = 1.03580745681136
$ws.Range("E16").Value = 1.034408103944936
$ws.Range("F16").Value = 1.041663695349706
$ws.Range("I16").Value = 1.036219922858047
$ws.Range("J16").Value = 1.031958027999996
$ws.Range("K16").Value = 1.039245655745946
$ws.Range("L16").Value = 1.037851299614842
$ws.Range("M16").Value = 1.045081143695719
$ws.Range("N16").Value = 1.014672243613843

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025902602866647
$ws.Range("D17").Value = 1.036013150799056
$ws.Range("E17").Value = 1.034716878174251
$ws.Range("F17").Value = 1.042035045621998
$ws.Range("I17").Value = 1.036291249532119
$ws.Range("J17").Value = 1.032152796450176
$ws.Range("K17").Value = 1.039378945286483
$ws.Range("L17").Value = 1.03808719630322
$ws.Range("M17").Value = 1.045379987731725
$ws.Range("N17").Value = 1.014736921193731

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026097360660425
$ws.Range("D18").Value = 1.036133101891511
$ws.Range("E18").Value = 1.034897077190655
$ws.Range("F18").Value = 1.042251784275559
$ws.Range("I18").Value = 1.036332680407585
$ws.Range("J18").Value = 1.03226637134941
$ws.Range("K18").Value = 1.039456578586186
$ws.Range("L18").Value = 1.038224808889731
$ws.Range("M18").Value = 1.045554360859387
$ws.Range("N18").Value = 1.014774633929296

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026163775816528
$ws.Range("D19").Value = 1.036173997519671
$ws.Range("E19").Value = 1.034958536672085
$ws.Range("F19").Value = 1.042325709628839
$ws.Range("I19").Value = 1.036346777934331
$ws.Range("J19").Value = 1.032305092291093
$ws.Range("K19").Value = 1.039483030454874
$ws.Range("L19").Value = 1.038271734246714
$ws.Range("M19").Value = 1.045613828140059
$ws.Range("N19").Value = 1.014787490847547

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025866782338362
$ws.Range("D20").Value = 1.035991084531346
$ws.Range("E20").Value = 1.034683739656218
$ws.Range("F20").Value = 1.041995189140034
$ws.Range("I20").Value = 1.036283614718716
$ws.Range("J20").Value = 1.032131902758728
$ws.Range("K20").Value = 1.039364656178838
$ws.Range("L20").Value = 1.038061884946019
$ws.Range("M20").Value = 1.045347918124276
$ws.Range("N20").Value = 1.014729983202108

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024901772934039
$ws.Range("D21").Value = 1.035396107030986
$ws.Range("E21").Value = 1.033791482668242
$ws.Range("F21").Value = 1.040922240951581
$ws.Range("I21").Value = 1.036076236415218
$ws.Range("J21").Value = 1.031568493932024
$ws.Range("K21").Value = 1.038978494056097
$ws.Range("L21").Value = 1.037379859454207
$ws.Range("M21").Value = 1.044484156255643
$ws.Range("N21").Value = 1.014542873082684

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024295458045913
$ws.Range("D22").Value = 1.035021793918589
$ws.Range("E22").Value = 1.033231363345005
$ws.Range("F22").Value = 1.040248875103393
$ws.Range("I22").Value = 1.035944299065164
$ws.Range("J22").Value = 1.031213991732669
$ws.Range("K22").Value = 1.038734697286228
$ws.Range("L22").Value = 1.03695121419249
$ws.Range("M22").Value = 1.043941647413473
$ws.Range("N22").Value = 1.014425118787888

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02461683631829
$ws.Range("D23").Value = 1.035220245251115
$ws.Range("E23").Value = 1.033528209428119
$ws.Range("F23").Value = 1.040605721071914
$ws.Range("I23").Value = 1.036014388214948
$ws.Range("J23").Value = 1.031401944771881
$ws.Range("K23").Value = 1.038864033005985
$ws.Range("L23").Value = 1.037178430502245
$ws.Range("M23").Value = 1.044229186868685
$ws.Range("N23").Value = 1.01448755292845

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025882967950849
$ws.Range("D24").Value = 1.036001055410734
$ws.Range("E24").Value = 1.03469871322893
$ws.Range("F24").Value = 1.042013198147594
$ws.Range("I24").Value = 1.036287065096172
$ws.Range("J24").Value = 1.032141343812099
$ws.Range("K24").Value = 1.03937111315847
$ws.Range("L24").Value = 1.038073322003224
$ws.Range("M24").Value = 1.045362408806633
$ws.Range("N24").Value = 1.014733118220864

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027354047476105
$ws.Range("D25").Value = 1.036906086296807
$ws.Range("E25").Value = 1.036060786881886
$ws.Range("F25").Value = 1.043651838737917
$ws.Range("I25").Value = 1.03659671091345
$ws.Range("J25").Value = 1.032998184784613
$ws.Range("K25").Value = 1.039955143518099
$ws.Range("L25").Value = 1.039112496811935
$ws.Range("M25").Value = 1.046679894671931
$ws.Range("N25").Value = 1.015017587218659
